# Auto-generated edit script: updates numeric cells H:N across multiple sheets
# per the authoritative diff (Excel COM interop).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 1356
$ws.Cells.Item(4, 9).Value = 1513
$ws.Cells.Item(4, 10).Value = 989.6667
$ws.Cells.Item(4, 11).Value = 1513
$ws.Cells.Item(4, 12).Value = 989.6667
$ws.Cells.Item(4, 13).Value = -1399
$ws.Cells.Item(4, 14).Value = -1217.6667
$ws.Cells.Item(12, 8).Value = 104.588234
$ws.Cells.Item(12, 9).Value = 104.588234
$ws.Cells.Item(12, 11).Value = 104.588234
$ws.Cells.Item(12, 13).Value = 65.411766
$ws.Cells.Item(15, 8).Value = 1988.5454
$ws.Cells.Item(15, 9).Value = 1988.5454
$ws.Cells.Item(15, 11).Value = 5965.6362
$ws.Cells.Item(15, 13).Value = -5796.6362
$ws.Cells.Item(48, 8).Value = 5019.016
$ws.Cells.Item(48, 10).Value = 5019.016
$ws.Cells.Item(48, 12).Value = 15057.048
$ws.Cells.Item(48, 14).Value = -15641.048
$ws.Cells.Item(56, 8).Value = 5019.016
$ws.Cells.Item(56, 10).Value = 5019.016
$ws.Cells.Item(56, 12).Value = 15057.048
$ws.Cells.Item(56, 14).Value = -16125.048
$ws.Cells.Item(76, 8).Value = 4630.7827
$ws.Cells.Item(76, 9).Value = 3789.75
$ws.Cells.Item(76, 10).Value = 5548.273
$ws.Cells.Item(76, 11).Value = 3789.75
$ws.Cells.Item(76, 12).Value = 5548.273
$ws.Cells.Item(76, 13).Value = -3474.75
$ws.Cells.Item(76, 14).Value = -6178.273
$ws.Cells.Item(79, 8).Value = 4630.7827
$ws.Cells.Item(79, 9).Value = 3789.75
$ws.Cells.Item(79, 10).Value = 5548.273
$ws.Cells.Item(79, 11).Value = 3789.75
$ws.Cells.Item(79, 12).Value = 5548.273
$ws.Cells.Item(79, 13).Value = -2697.75
$ws.Cells.Item(79, 14).Value = -7732.273
$ws.Cells.Item(106, 8).Value = 7275.2593
$ws.Cells.Item(106, 9).Value = 6978.154
$ws.Cells.Item(106, 10).Value = 15000
$ws.Cells.Item(106, 11).Value = 6978.154
$ws.Cells.Item(106, 12).Value = 15000
$ws.Cells.Item(106, 13).Value = -6347.154
$ws.Cells.Item(106, 14).Value = -16262
$ws.Cells.Item(109, 8).Value = 205287760
$ws.Cells.Item(109, 10).Value = 205287760
$ws.Cells.Item(109, 12).Value = 205287760
$ws.Cells.Item(109, 14).Value = -205290534
$ws.Cells.Item(132, 8).Value = 3133.463
$ws.Cells.Item(132, 9).Value = 2983.9575
$ws.Cells.Item(132, 11).Value = 8951.872499999999
$ws.Cells.Item(132, 13).Value = -6421.872499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 4884.95
$ws.Cells.Item(45, 9).Value = 3952.3333
$ws.Cells.Item(45, 11).Value = 3952.3333
$ws.Cells.Item(45, 13).Value = -3575.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3227.9546
$ws.Cells.Item(20, 9).Value = 1794.8462
$ws.Cells.Item(20, 10).Value = 5298
$ws.Cells.Item(20, 11).Value = 1794.8462
$ws.Cells.Item(20, 12).Value = 5298
$ws.Cells.Item(20, 13).Value = -1547.8462
$ws.Cells.Item(20, 14).Value = -5792
$ws.Cells.Item(46, 8).Value = 8250
$ws.Cells.Item(46, 10).Value = 8250
$ws.Cells.Item(46, 12).Value = 8250
$ws.Cells.Item(46, 14).Value = -8846
$ws.Cells.Item(82, 8).Value = 22155.125
$ws.Cells.Item(82, 10).Value = 59965.5
$ws.Cells.Item(82, 12).Value = 59965.5
$ws.Cells.Item(82, 14).Value = -60731.5
$ws.Cells.Item(85, 8).Value = 22155.125
$ws.Cells.Item(85, 10).Value = 59965.5
$ws.Cells.Item(85, 12).Value = 59965.5
$ws.Cells.Item(85, 14).Value = -62617.5
$ws.Cells.Item(134, 8).Value = 11367.588
$ws.Cells.Item(134, 9).Value = 12635
$ws.Cells.Item(134, 11).Value = 37905
$ws.Cells.Item(134, 13).Value = -35370

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5390.3687
$ws.Cells.Item(31, 9).Value = 5446.2583
$ws.Cells.Item(31, 11).Value = 5446.2583
$ws.Cells.Item(31, 13).Value = -5151.2583
$ws.Cells.Item(34, 8).Value = 5390.3687
$ws.Cells.Item(34, 9).Value = 5446.2583
$ws.Cells.Item(34, 11).Value = 5446.2583
$ws.Cells.Item(34, 13).Value = -5244.2583
$ws.Cells.Item(134, 8).Value = 13756.3
$ws.Cells.Item(134, 9).Value = 18937.572
$ws.Cells.Item(134, 11).Value = 56812.716
$ws.Cells.Item(134, 13).Value = -54277.716

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(56, 8).Value = 6606.8696
$ws.Cells.Item(56, 9).Value = 6606.8696
$ws.Cells.Item(56, 11).Value = 6606.8696
$ws.Cells.Item(56, 13).Value = -6076.8696

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 9803.951999999999
$ws.Cells.Item(70, 9).Value = 9487.625
$ws.Cells.Item(70, 10).Value = 9998.615
$ws.Cells.Item(70, 11).Value = 9487.625
$ws.Cells.Item(70, 12).Value = 9998.615
$ws.Cells.Item(70, 13).Value = -9217.625
$ws.Cells.Item(70, 14).Value = -10538.615
$ws.Cells.Item(73, 8).Value = 9803.951999999999
$ws.Cells.Item(73, 9).Value = 9487.625
$ws.Cells.Item(73, 10).Value = 9998.615
$ws.Cells.Item(73, 11).Value = 9487.625
$ws.Cells.Item(73, 12).Value = 9998.615
$ws.Cells.Item(73, 13).Value = -8551.625
$ws.Cells.Item(73, 14).Value = -11870.615
$ws.Cells.Item(102, 8).Value = 5863.407
$ws.Cells.Item(102, 9).Value = 6838
$ws.Cells.Item(102, 11).Value = 6838
$ws.Cells.Item(102, 13).Value = -5216
$ws.Cells.Item(122, 8).Value = 10216.381
$ws.Cells.Item(122, 9).Value = 6660.2104
$ws.Cells.Item(122, 10).Value = 44000
$ws.Cells.Item(122, 11).Value = 19980.6312
$ws.Cells.Item(122, 12).Value = 132000
$ws.Cells.Item(122, 13).Value = -17530.6312
$ws.Cells.Item(122, 14).Value = -136900
$ws.Cells.Item(126, 8).Value = 9413.444
$ws.Cells.Item(126, 10).Value = 3301.182
$ws.Cells.Item(126, 12).Value = 9903.545999999998
$ws.Cells.Item(126, 14).Value = -14843.546
$ws.Cells.Item(132, 8).Value = 5284.609
$ws.Cells.Item(132, 9).Value = 5311.7144
$ws.Cells.Item(132, 11).Value = 15935.1432
$ws.Cells.Item(132, 13).Value = -13405.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 18765.172
$ws.Cells.Item(7, 9).Value = 22866.38
$ws.Cells.Item(7, 10).Value = 7999.5
$ws.Cells.Item(7, 11).Value = 22866.38
$ws.Cells.Item(7, 12).Value = 7999.5
$ws.Cells.Item(7, 13).Value = -22754.38
$ws.Cells.Item(7, 14).Value = -8223.5
$ws.Cells.Item(16, 8).Value = 1275.9
$ws.Cells.Item(16, 9).Value = 1308.5416
$ws.Cells.Item(16, 10).Value = 1145.3334
$ws.Cells.Item(16, 11).Value = 1308.5416
$ws.Cells.Item(16, 12).Value = 1145.3334
$ws.Cells.Item(16, 13).Value = -1138.5416
$ws.Cells.Item(16, 14).Value = -1485.3334
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 14).ClearContents()
$ws.Cells.Item(43, 8).Value = 32600
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 14).ClearContents()
$ws.Cells.Item(55, 8).Value = 866.25
$ws.Cells.Item(55, 9).Value = 229.88235
$ws.Cells.Item(55, 10).Value = 2411.7144
$ws.Cells.Item(55, 11).Value = 229.88235
$ws.Cells.Item(55, 12).Value = 2411.7144
$ws.Cells.Item(55, 13).Value = -56.88235
$ws.Cells.Item(55, 14).Value = -2757.7144
$ws.Cells.Item(87, 8).Value = 33999
$ws.Cells.Item(87, 9).Value = 33999
$ws.Cells.Item(87, 11).Value = 33999
$ws.Cells.Item(87, 13).Value = -32876
$ws.Cells.Item(90, 8).Value = 33999
$ws.Cells.Item(90, 9).Value = 33999
$ws.Cells.Item(90, 11).Value = 101997
$ws.Cells.Item(90, 13).Value = -96381
$ws.Cells.Item(126, 8).Value = 18765.172
$ws.Cells.Item(126, 9).Value = 22866.38
$ws.Cells.Item(126, 10).Value = 7999.5
$ws.Cells.Item(126, 11).Value = 68599.14
$ws.Cells.Item(126, 12).Value = 23998.5
$ws.Cells.Item(126, 13).Value = -66129.14
$ws.Cells.Item(126, 14).Value = -28938.5
$ws.Cells.Item(132, 8).Value = 1658744.8
$ws.Cells.Item(132, 9).Value = 3725238.2
$ws.Cells.Item(132, 10).Value = 5550
$ws.Cells.Item(132, 11).Value = 11175714.6
$ws.Cells.Item(132, 12).Value = 16650
$ws.Cells.Item(132, 13).Value = -11173184.6
$ws.Cells.Item(132, 14).Value = -21710
$ws.Cells.Item(136, 8).Value = 7067.9287
$ws.Cells.Item(136, 9).Value = 7624.5
$ws.Cells.Item(136, 11).Value = 22873.5
$ws.Cells.Item(136, 13).Value = -20323.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 14).ClearContents()
$ws.Cells.Item(58, 8).Value = 6813197
$ws.Cells.Item(58, 9).Value = 18492.5
$ws.Cells.Item(58, 10).Value = 11343000
$ws.Cells.Item(58, 11).Value = 18492.5
$ws.Cells.Item(58, 12).Value = 11343000
$ws.Cells.Item(58, 13).Value = -18184.5
$ws.Cells.Item(58, 14).Value = -11343616
$ws.Cells.Item(132, 8).Value = 14825.533
$ws.Cells.Item(132, 9).Value = 24526.5
$ws.Cells.Item(132, 10).Value = 6337.1875
$ws.Cells.Item(132, 11).Value = 73579.5
$ws.Cells.Item(132, 12).Value = 19011.5625
$ws.Cells.Item(132, 13).Value = -71049.5
$ws.Cells.Item(132, 14).Value = -24071.5625
$ws.Cells.Item(135, 8).Value = 51200
$ws.Cells.Item(135, 10).Value = 51200
$ws.Cells.Item(135, 12).Value = 51200
$ws.Cells.Item(135, 14).Value = -61340
